$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings
# (e.g. "0.997", "10.01") are stored as text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '53.994.43'
$ws.Range("E2").Value = '  +1.26%  '
$ws.Range("D3").Value = '2.279.03'
$ws.Range("E3").Value = '  +3.55%  '
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '495.27'
$ws.Range("E5").Value = '  +3.15%  '
$ws.Range("D6").Value = '127.63'
$ws.Range("E6").Value = '  +2.67%  '
$ws.Range("D7").Value = '0.994'
$ws.Range("E7").Value = '  -0.35%  '
$ws.Range("E8").Value = '  +2.52%  '
$ws.Range("D9").Value = '2.276.64'
$ws.Range("E9").Value = '  +3.09%  '
$ws.Range("D10").Value = '0.0951'
$ws.Range("E10").Value = '  +4.64%  '
$ws.Range("E11").Value = '  +2.43%  '
$ws.Range("D12").Value = '0.324'
$ws.Range("E12").Value = '  +4.19%  '
$ws.Range("D13").Value = '4.65'
$ws.Range("E13").Value = '  +0.05%  '
$ws.Range("D14").Value = '2.667.02'
$ws.Range("E14").Value = '  +2.91%  '
$ws.Range("D15").Value = '21.74'
$ws.Range("E15").Value = '  +3.98%  '
$ws.Range("D16").Value = '53.961.75'
$ws.Range("E16").Value = '  +1.54%  '
$ws.Range("E17").Value = '  +1.95%  '
$ws.Range("D18").Value = '2.267.06'
$ws.Range("E18").Value = '  +2.68%  '
$ws.Range("D19").Value = '10.01'
$ws.Range("E19").Value = '  +5.71%  '
$ws.Range("D20").Value = '4.09'
$ws.Range("E20").Value = '  +4.45%  '
$ws.Range("D21").Value = '300.41'
$ws.Range("E21").Value = '  +1.71%  '
$ws.Range("D22").Value = '6.42'
$ws.Range("E22").Value = '  +6.35%  '
$ws.Range("D23").Value = '0.997'
$ws.Range("E23").Value = '  -0.31%  '
$ws.Range("D24").Value = '5.39'
$ws.Range("E24").Value = '  -2.42%  '
$ws.Range("D25").Value = '62.16'
$ws.Range("E25").Value = '  -0.74%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("E27").Value = '  +2.89%  '
$ws.Range("D28").Value = '2.378.04'
$ws.Range("E28").Value = '  +3.31%  '
$ws.Range("E29").Value = '  +4.81%  '
$ws.Range("D30").Value = '7.06'
$ws.Range("E30").Value = '  +2.08%  '
$ws.Range("D31").Value = '167.97'
$ws.Range("E31").Value = '  +0.63%  '
$ws.Range("E32").Value = '  +2.50%  '
$ws.Range("D33").Value = '0.0₃0687'
$ws.Range("E33").Value = '  +4.13%  '
$ws.Range("D34").Value = '5.87'
$ws.Range("E34").Value = '  +4.38%  '
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("D36").Value = '0.996'
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("E37").Value = '  +3.39%  '
$ws.Range("E38").Value = '  +3.24%  '
$ws.Range("D39").Value = '0.912'
$ws.Range("E39").Value = '  +11.72%  '
$ws.Range("E40").Value = '  +4.72%  '
$ws.Range("E41").Value = '  +4.97%  '
$ws.Range("D42").Value = '35.61'
$ws.Range("E42").Value = '  +0.23%  '
$ws.Range("E43").Value = '  +4.38%  '
$ws.Range("E44").Value = '  +2.85%  '
$ws.Range("D45").Value = '3.35'
$ws.Range("E45").Value = '  +3.77%  '
$ws.Range("D46").Value = '125.93'
$ws.Range("E46").Value = '  +3.84%  '
$ws.Range("D47").Value = '4.76'
$ws.Range("E47").Value = '  +1.46%  '
$ws.Range("E48").Value = '  +1.57%  '
$ws.Range("E49").Value = '  +3.23%  '
$ws.Range("D50").Value = '237.40'
$ws.Range("E50").Value = '  +5.03%  '
$ws.Range("E51").Value = '  +3.81%  '

# Restore column D to the default ("Normal") style so no stray
# number-format style is left on the cells.
$ws.Range("D2:D51").Style = "Normal"
